$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (E1:I1), matching the bold/bordered header style used by A1:D1
$ws.Range("E1").Value = "Src_DeliveryReceipt"
$ws.Range("F1").Value = "Src_Remaining"
$ws.Range("G1").Value = "Src_Transfers"
$ws.Range("H1").Value = "Src_Beverages"
$ws.Range("I1").Value = "DR Price"

# Copy the header style (bold, border, centered) from D1 onto the new header cells
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:I1").PasteSpecial(-4122) | Out-Null

# New data cells for row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
